# Insert a new weekly record at row 143 of the "Hortaliza, Vega Monumental
# Concepción - Zapallo italiano" sheet. Inserting the row shifts the
# existing rows 143-180 down to 144-181 (formatting - e.g. the date style
# on column D - comes along automatically with the native Insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(143).Insert()

$ws.Cells.Item(143, 1).Value = 11
$ws.Cells.Item(143, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(143, 3).Value = "Bíobío"
$ws.Cells.Item(143, 4).Value = 44932
$ws.Cells.Item(143, 5).Value = 8
$ws.Cells.Item(143, 6).Value = 100112032
$ws.Cells.Item(143, 7).Value = "Zapallo italiano"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 450
$ws.Cells.Item(143, 11).Value = 10000
$ws.Cells.Item(143, 12).Value = 11000
$ws.Cells.Item(143, 13).Value = 10556
$ws.Cells.Item(143, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(143, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(143, 16).Value = 176
$ws.Cells.Item(143, 17).Value = 60
$ws.Cells.Item(143, 18).Value = "Hortaliza"
